# The deck currently uses the "Integral" design (ppt/theme/theme2.xml, the
# theme bound to the slide master / presentation) while ppt/theme/theme1.xml
# still holds the stock default "Office Theme" palette (only wired to the
# notes master). The edit swaps the two: the design actually applied to the
# slides switches from "Integral" to the default "Office Theme" color
# palette. We drive this by rewriting the theme color scheme that is bound
# to the slides (reachable via Slide.ThemeColorScheme / the slide master's
# ColorScheme - they are the same underlying theme), setting every slot to
# the stock Office theme's RGB values.

function Set-ThemeColor {
    param($ColorScheme, [int]$Index, [string]$HexRGB)

    $r = [Convert]::ToInt32($HexRGB.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($HexRGB.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($HexRGB.Substring(4, 2), 16)

    # PowerPoint's RGB color values are packed as 0xBBGGRR (same convention
    # as the classic VBA RGB() function), i.e. low byte = red.
    $ColorScheme.Item($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Standard Office theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeColor $tcs 1  "000000"
Set-ThemeColor $tcs 2  "FFFFFF"
Set-ThemeColor $tcs 3  "44546A"
Set-ThemeColor $tcs 4  "E7E6E6"
Set-ThemeColor $tcs 5  "5B9BD5"
Set-ThemeColor $tcs 6  "ED7D31"
Set-ThemeColor $tcs 7  "A5A5A5"
Set-ThemeColor $tcs 8  "FFC000"
Set-ThemeColor $tcs 9  "4472C4"
Set-ThemeColor $tcs 10 "70AD47"
Set-ThemeColor $tcs 11 "0563C1"
Set-ThemeColor $tcs 12 "954F72"
